# Update cryptocurrency price/volume figures to the latest scraped values.
# Cells in columns D (Price) and E (Volume(1h)) are stored as text, so we
# force a text number format before writing the new value to avoid Excel
# auto-converting the look-alike numeric/percentage strings into numbers.
function Set-TextValue {
    param($ws, $ref, $val)
    $c = $ws.Range($ref)
    $c.NumberFormat = "@"
    $c.Value = $val
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

Set-TextValue $ws "D2" "273.11"
Set-TextValue $ws "D3" "26.89"
Set-TextValue $ws "E3" "-0.96%"
Set-TextValue $ws "D4" "4.727"
Set-TextValue $ws "E4" "0.43%"
Set-TextValue $ws "D5" "0.06187"
Set-TextValue $ws "E5" "-0.41%"
Set-TextValue $ws "E6" "0.54%"
Set-TextValue $ws "D7" "0.8630"
Set-TextValue $ws "E7" "1.50%"
Set-TextValue $ws "D8" "0.9122"
Set-TextValue $ws "E8" "0.54%"
Set-TextValue $ws "D9" "0.1440"
Set-TextValue $ws "E9" "2.77%"
Set-TextValue $ws "D10" "0.05322"
Set-TextValue $ws "E10" "12.88%"
Set-TextValue $ws "D11" "0.07159"
Set-TextValue $ws "E11" "0.98%"
Set-TextValue $ws "D12" "0.03180"
Set-TextValue $ws "E12" "0.12%"
Set-TextValue $ws "D13" "0.09054"
Set-TextValue $ws "E13" "-0.09%"
Set-TextValue $ws "D14" "0.001530"
Set-TextValue $ws "E14" "-0.04%"
Set-TextValue $ws "D15" "0.0006073"
Set-TextValue $ws "E15" "-1.71%"
Set-TextValue $ws "D16" "0.005988"
Set-TextValue $ws "E16" "-0.05%"
Set-TextValue $ws "D17" "3.474"
Set-TextValue $ws "E17" "0.21%"
Set-TextValue $ws "D18" "3.192"
Set-TextValue $ws "E18" "0.63%"
Set-TextValue $ws "E19" "4.00%"
Set-TextValue $ws "E20" "-0.68%"
Set-TextValue $ws "E21" "1.35%"
Set-TextValue $ws "D22" "3.848"
Set-TextValue $ws "E22" "-5.95%"
Set-TextValue $ws "D23" "0.04252"
Set-TextValue $ws "E23" "0.21%"
Set-TextValue $ws "D24" "0.001176"
Set-TextValue $ws "E24" "-3.56%"
Set-TextValue $ws "D25" "0.004194"
Set-TextValue $ws "E25" "1.90%"
Set-TextValue $ws "E26" "-0.16%"
Set-TextValue $ws "D40" "0.03975"
Set-TextValue $ws "E40" "1.95%"
Set-TextValue $ws "D41" "0.006211"
Set-TextValue $ws "E41" "50.24%"
Set-TextValue $ws "D42" "0.1130"
Set-TextValue $ws "E42" "1.67%"
Set-TextValue $ws "D43" "0.002298"
Set-TextValue $ws "E43" "5.24%"
Set-TextValue $ws "D44" "0.01276"
Set-TextValue $ws "E44" "-4.85%"
Set-TextValue $ws "D45" "0.00005119"
Set-TextValue $ws "E45" "-1.06%"
Set-TextValue $ws "E46" "-0.17%"
Set-TextValue $ws "D47" "0.8973"
Set-TextValue $ws "E47" "450.58%"
Set-TextValue $ws "E48" "-14.69%"
Set-TextValue $ws "E49" "-0.17%"
Set-TextValue $ws "E50" "-0.17%"
